$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Dittmann")

# Update hours for row 11 (B11=16) from 1.5 to 3.5
$ws.Range("D11").Value = 3.5

# Fill in row 12 (B12=17) with the new AWS documentation task and its hours
$ws.Range("C12").Value = "Dokumentation AWS"
$ws.Range("D12").Value = 4

# Move the active selection to C13 (next empty row) as in the authored edit
$ws.Range("C13").Select()
